$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same duplicated rows that
# need their "想去人数" (want-to-go count) values refreshed.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 1492
    $ws.Range("F9").Value = 280
}
